$wb = $excel.ActiveWorkbook

# --- Rename headers on existing sheets ---
$ws1 = $wb.Worksheets.Item("Weekly Quantity")
$ws1.Range("B1").Value = "Weekly_PO_Qty"

$ws2 = $wb.Worksheets.Item("Monthly Trend")
$ws2.Range("B1").Value = "Monthly_PO_Qty"

# --- Add the new "PO Forecast" sheet after "Monthly Trend" (i.e. at the end) ---
$ws3 = $wb.Worksheets.Add($null, $ws2)
$ws3.Name = "PO Forecast"

# Match the sheet-level layout of the existing sheets
$ws3.Outline.SummaryRow = 1
$ws3.Outline.SummaryColumn = 1
$ws3.PageSetup.LeftMargin = 54
$ws3.PageSetup.RightMargin = 54
$ws3.PageSetup.TopMargin = 72
$ws3.PageSetup.BottomMargin = 72
$ws3.PageSetup.HeaderMargin = 36
$ws3.PageSetup.FooterMargin = 36

# Reuse the bold/centered header style from the existing sheets for row 1
$ws1.Range("A1:B1").Copy()
$ws3.Range("A1:D1").PasteSpecial(-4122)

# Reuse the date-formatted style from column A of the existing sheets
$ws1.Range("A2").Copy()
$ws3.Range("A2:A12").PasteSpecial(-4122)

# --- Header row ---
$ws3.Range("A1").Value = "ds"
$ws3.Range("B1").Value = "PO_Forecast"
$ws3.Range("C1").Value = "yhat_lower"
$ws3.Range("D1").Value = "yhat_upper"

# --- Data rows ---
$ws3.Range("A2").Value = 45578.99999999999
$ws3.Range("B2").Value = 244
$ws3.Range("C2").Value = 190.8344703859924
$ws3.Range("D2").Value = 296.8280195545025

$ws3.Range("A3").Value = 45592.99999999999
$ws3.Range("B3").Value = 71
$ws3.Range("C3").Value = 17.23528348013812
$ws3.Range("D3").Value = 124.2822442046863

$ws3.Range("A4").Value = 45599.99999999999
$ws3.Range("B4").Value = 0
$ws3.Range("C4").Value = -69.37757921891719
$ws3.Range("D4").Value = 37.90555896069741

$ws3.Range("A5").Value = 45606.99999999999
$ws3.Range("B5").Value = 0
$ws3.Range("C5").Value = -158.3659519967824
$ws3.Range("D5").Value = -50.71965738237075

$ws3.Range("A6").Value = 45613.99999999999
$ws3.Range("B6").Value = 0
$ws3.Range("C6").Value = -244.4650351062514
$ws3.Range("D6").Value = -137.4284393150575

$ws3.Range("A7").Value = 45620.99999999999
$ws3.Range("B7").Value = 0
$ws3.Range("C7").Value = -327.3788401634081
$ws3.Range("D7").Value = -220.7320410304109

$ws3.Range("A8").Value = 45627.99999999999
$ws3.Range("B8").Value = 0
$ws3.Range("C8").Value = -413.9181235107314
$ws3.Range("D8").Value = -305.6931337396678

$ws3.Range("A9").Value = 45634.99999999999
$ws3.Range("B9").Value = 0
$ws3.Range("C9").Value = -499.3987300985029
$ws3.Range("D9").Value = -395.5155311491887

$ws3.Range("A10").Value = 45641.99999999999
$ws3.Range("B10").Value = 0
$ws3.Range("C10").Value = -589.4978927676215
$ws3.Range("D10").Value = -482.5280066577504

$ws3.Range("A11").Value = 45648.99999999999
$ws3.Range("B11").Value = 0
$ws3.Range("C11").Value = -670.8478529374768
$ws3.Range("D11").Value = -568.667435996218

$ws3.Range("A12").Value = 45655.99999999999
$ws3.Range("B12").Value = 0
$ws3.Range("C12").Value = -755.7271149400029
$ws3.Range("D12").Value = -646.3827390028505

# Leave the selection on A1 of the new sheet, like a freshly added sheet
[void]$ws3.Range("A1").Select()
